$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: DocNum, Name, Period, ValorMora, SalarioBasico
$rows = @(
    @("1002319881", "KEIVER BOSSIO BALLESTEROS", "2212", 34666, 1000000),
    @("1050952836", "CINDY MILENA LARA ESPITALETA", "2303", 16000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2304", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2305", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2306", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2307", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2308", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2309", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2310", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2311", 60000, 1500000),
    @("1044928283", "GENESIS TORRES RICO", "2312", 38000, 1500000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 3).Value = $data[0]
    $ws.Cells.Item($r, 4).Value = $data[1]
    $ws.Cells.Item($r, 5).Value = $data[2]
    $ws.Cells.Item($r, 6).Value = $data[3]
    $ws.Cells.Item($r, 7).Value = $data[4]
}
